$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at the top; this shifts all existing rows down by 1
$ws.Rows("1:1").Insert()

# Copy the header-row formatting (now on row 2) onto the new row 1
$ws.Range("A2:L2").Copy()
$ws.Range("A1:L1").PasteSpecial(-4122)

# Fill new row 1 with sequential numbers 0-11
for ($i = 0; $i -lt 12; $i++) {
    $ws.Cells.Item(1, $i + 1).Value = $i
}

# The old header row (now row 2) loses its bold/border style
$ws.Range("A2:L2").Style = "Normal"

# Clear the thread_size / material_surface labels from the (now) row 2
$ws.Range("K2").ClearContents()
$ws.Range("L2").ClearContents()
